$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1.xml)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 5390.387531524985
$ws.Range("B2").Value = 6663.107410080002
$ws.Range("E2").Value = 14816.31105402212
$ws.Range("I2").Value = 36846.71726986
$ws.Range("M2").Value = 11706.203357365
$ws.Range("N2").Value = 3919.172050548081
$ws.Range("O2").Value = 6924.712879078098

# Sheet "2030" (sheet2.xml)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 5390.387531524985
$ws.Range("B2").Value = 8293.296836026202
$ws.Range("E2").Value = 40066.17459638815
$ws.Range("I2").Value = 67226.01974586057
$ws.Range("M2").Value = 22423.598584156
$ws.Range("N2").Value = 10156.34074110075
$ws.Range("O2").Value = 12108.04910209266

# Sheet "2035" (sheet3.xml)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 5390.387531524985
$ws.Range("B2").Value = 8293.296836026202
$ws.Range("E2").Value = 60397.03343082713
$ws.Range("G2").Value = 7864.0611328728
$ws.Range("I2").Value = 87485.02868918961
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 28551.55208405499
$ws.Range("N2").Value = 12959.04701199569
$ws.Range("O2").Value = 15398.65911648144

# Sheet "2040" (sheet4.xml)
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 5390.387531524985
$ws.Range("B2").Value = 8293.296836026202
$ws.Range("E2").Value = 60397.03343082713
$ws.Range("G2").Value = 7864.0611328728
$ws.Range("I2").Value = 87485.02868918961
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 28551.55208405499
$ws.Range("N2").Value = 12959.04701199569
$ws.Range("O2").Value = 15398.65911648144

# Sheet "2045" (sheet5.xml)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 5390.387531524985
$ws.Range("B2").Value = 8293.296836026202
$ws.Range("E2").Value = 60397.03343082713
$ws.Range("G2").Value = 7864.0611328728
$ws.Range("I2").Value = 87485.02868918961
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 28551.55208405499
$ws.Range("N2").Value = 12959.04701199569
$ws.Range("O2").Value = 15398.65911648144

# Sheet "2050" (sheet6.xml)
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 5390.387531524985
$ws.Range("B2").Value = 8293.296836026202
$ws.Range("E2").Value = 60397.03343082713
$ws.Range("G2").Value = 7864.0611328728
$ws.Range("I2").Value = 87485.02868918961
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 28551.55208405499
$ws.Range("N2").Value = 12959.04701199569
$ws.Range("O2").Value = 15398.65911648144
